$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.161.67'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '3.925.99'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '486.10'
$ws.Range('E5').Value = '  +3.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.43'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.730'
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('E10').Value = '  +3.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000362'
$ws.Range('E11').Value = '  +6.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.78'
$ws.Range('E12').Value = '  -0.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.67'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Value = '4.552.76'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.94'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = '3.925.71'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.15'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').Value = '68.328.43'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '447.85'
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.76'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.53'
$ws.Range('E25').Value = '  +14.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.62'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('E27').Value = '  +11.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.05'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.51'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.131'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '690.20'
$ws.Range('E32').Value = '  -6.41%  '
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').Value = '0.0₃0951'
$ws.Range('E34').Value = '  +20.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.89'
$ws.Range('E35').Value = '  -3.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.86'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.74'
$ws.Range('E37').Value = '  +6.23%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.150'
$ws.Range('E38').Value = '  -4.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('E41').Value = '  +11.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.83'
$ws.Range('E42').Value = '  +13.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.10'
$ws.Range('E43').Value = '  -4.92%  '
$ws.Range('E44').Value = '  +5.91%  '
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.44'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.15'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '146.09'
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.16'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.84'
$ws.Range('E51').Value = '  -2.20%  '
